# Adding messages for svat_unbalanced and svat_unbalanced_nok
#
# Inserts two new rows into the i18n table (right before the existing
# "svat_ok_activity_start" row) holding the keys/pt-text for two new
# translation messages, then grows the table + dimension to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 542 (pushes everything from 542 down to 544+)
$ws.Rows("542:543").Insert()

# Column A holds the i18n keys, column B holds the Portuguese text.
# Write all of column A first, then column B, so the new shared-string
# entries are appended in key/key/value/value order.
$ws.Range("A542").Value = "svat_unbalanced"
$ws.Range("A543").Value = "svat_unbalanced_nok"
$ws.Range("B542").Value = "Teste a movimentos finalizados em moeda estrangeira"
$ws.Range("B543").Value = "Verificamos que alguns movimentos efectuados em moeda estrangeira não estão consistentes, por favor reabra e corrija os seguintes movimentos:"

# Row heights for the two new wrapped-text rows
$ws.Rows(542).RowHeight = 34
$ws.Rows(543).RowHeight = 51

# Grow the "i18n" table so the new rows are included in it
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:G582"))

# Leave the selection on the last edited cell
$ws.Range("B543").Select()
